# Applies the "1_3_PSP_Sheet_이상우.xlsx" update:
#  - corrects the 10/25 -> 10/22 and 10/27 -> 10/24 log dates
#  - rewrites the activity descriptions for rows 22-23
#  - fixes the start/stop/interruption/delta numbers for rows 22-23
#  - fills in two brand-new log rows (24-25) that were blank before
#  - scrolls the sheet so row 7 is the first visible row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: 10월 22일 / 설문지 자료 정리 ------------------------------
$ws.Range("A22").Value2 = "10월 22일"
$ws.Range("B22").Value2 = 0.77083333333333337
$ws.Range("C22").Value2 = 0.875
$ws.Range("D22").Value2 = 20
$ws.Range("E22").Value2 = 150
$ws.Range("F22").Value2 = "설문지 자료 정리"

# --- Row 23: 10월 24일 / 설문지 자료 토대로 Initial Data set 작성 -----
$ws.Range("A23").Value2 = "10월 24일"
$ws.Range("B23").Value2 = 0.45833333333333331
$ws.Range("C23").Value2 = 0.78402777777777777
$ws.Range("D23").Value2 = 70
$ws.Range("E23").Value2 = 469
$ws.Range("F23").Value2 = "설문지 자료 토대로 Initial Data set 작성"

# --- Row 24 (previously blank): 10월 29일 / web1 4강 ------------------
$ws.Range("A24").Value2 = "10월 29일"
$ws.Range("B24").Value2 = 0.70833333333333337
$ws.Range("C24").Value2 = 0.79166666666666663
$ws.Range("D24").Value2 = 0
$ws.Range("E24").Value2 = 120
$ws.Range("F24").Value2 = "web1 4강"

# --- Row 25 (previously blank): 11월 3일 / web1 6강 -------------------
$ws.Range("A25").Value2 = "11월 3일"
$ws.Range("B25").Value2 = 0.5
$ws.Range("C25").Value2 = 0.16388888888888889
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 236
$ws.Range("F25").Value2 = "web1 6강"

# --- Scroll the window so row 7 becomes the top-left visible cell -----
$ws.Range("A7").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

# restore the original selection
$ws.Range("F26").Select()
